$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix A92: the date/time stamp had a stray time-of-day component; ---
# --- it should just be the plain date-stamp used throughout the sheet. ---
$ws.Range("A92").Value = 45469.2916666667

# --- Append a new row (93) with the latest scraped data point ---

# Copy the date/time number format from A92 onto A93 first so the new
# cell re-uses the existing "yyyy-mm-dd hh:mm:ss" style rather than a
# freshly minted one.
$ws.Range("A92").Copy() | Out-Null
$ws.Range("A93").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A93").Value = 45470.2917476852
$ws.Range("B93").Value = 600
$ws.Range("C93").Value = 6.26000022888184
$ws.Range("D93").Value = 6.26000022888184
$ws.Range("E93").Value = 6.26000022888184
$ws.Range("F93").Value = 6.26000022888184

# G93/H93 hold text (shared-string) values, matching the rest of the
# sheet, even though G93's content looks numeric.
$ws.Range("G93").NumberFormat = "@"
$ws.Range("G93").Value = "6.26000022888184"
$ws.Range("G93").Style = $ws.Range("B92").Style

$ws.Range("H93").Value = "PAL.MI"
